$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Value = "Datatype aaa"
$ws.Range("B4").Value = "BigDecimal"
$ws.Range("C4").Value = "aaaaa"
$ws.Range("B3:C4").Borders.LineStyle = 1
$ws.Range("B3:C3").Merge()
